$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table runs from row 2 through row 120 (columns A:J), one row per
# day. Extend it with six more daily rows (121-126), repeating the last
# row's B:J readings and continuing the date sequence in column A.

$srcRow = 120
$numNewRows = 6

# Capture the source row's values (Value2 avoids locale/date-string coercion).
$srcValues = @()
for ($col = 1; $col -le 10; $col++) {
    $srcValues += ,$ws.Cells.Item($srcRow, $col).Value2
}
$startDate = $srcValues[0] + 1

for ($i = 0; $i -lt $numNewRows; $i++) {
    $dstRow = $srcRow + 1 + $i

    for ($col = 1; $col -le 10; $col++) {
        if ($col -eq 1) {
            $ws.Cells.Item($dstRow, $col).Value = $startDate + $i
        } else {
            $ws.Cells.Item($dstRow, $col).Value = $srcValues[$col - 1]
        }
    }

    # Copy the source row's formatting (number format/style) onto the new
    # row without touching the values we just wrote.
    $ws.Range("A$srcRow`:J$srcRow").Copy() | Out-Null
    $ws.Range("A$dstRow").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
